# downloadTotalForm.xlsx edit:
#  - Insert two new "yes/no" header fields (PinCode, AddressState) into the
#    header row, pushing the existing headers from Pan Card onward two
#    columns to the right (AA1:AS1 -> AC1:AU1), and inserting the new
#    headers at AA1/AB1.
#  - Append two new trailing "yes" answer columns (AT/AU) to every data row.
#  - Fix the "Cash" payment-mode row (row 14) so the previously "no"
#    answers in columns B, G, H, K become "yes" (collector-name auth no
#    longer required for cash transactions).
#  - Update the sheet's saved selection to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: shift Pan Card..Cheque/DD photo2 two columns to the
#    right (AA1:AS1 -> AC1:AU1), then write the two brand-new headers
#    into the freed-up AA1/AB1 cells.
# ---------------------------------------------------------------------

$oldHeaders = @(
    "Pan Card [ yes,no]",
    "Pan Card Remark [ yes,no ]",
    "Amount       [ yes,no ]",
    "Amount in Words [ yes,no ]",
    "Collector Name                [ yes,no ]",
    "Collector Phone                [ yes,no ]",
    "Nature of Donation                 [ yes,no ]",
    "Party Unit [ yes,no ]",
    "Location [ yes,no ]",
    "Payment realize date [ yes,no ]",
    "Receipt Number [ yes,no ]",
    "Transaction Valid [ yes,no ]",
    "Created By [ yes,no ]",
    "Created At [ yes,no ]",
    "Cheque Bounce Remark [yes,no ]",
    "Reverse Remark      [ yes,no ]",
    "Pan Card Photo [ yes,no ]",
    "Cheque/DD photo1 [ yes,no ]",
    "Cheque/DD photo2 [ yes,no ]"
)

function Get-ColumnName($colNum) {
    $dividend = $colNum
    $columnName = ""
    while ($dividend -gt 0) {
        $modulo = ($dividend - 1) % 26
        $columnName = [char](65 + $modulo) + $columnName
        $dividend = [int](($dividend - $modulo) / 26)
    }
    return $columnName
}

# New home for the shifted headers starts at column 29 (AC).
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = Get-ColumnName (29 + $i)
    $addr = $col + "1"
    $ws.Range("Z1").Copy($ws.Range($addr))
    $ws.Range($addr).Value = $oldHeaders[$i]
}

# New header cells, formatted like the rest of row 1 (bold + wrap).
$ws.Range("Z1").Copy($ws.Range("AA1"))
$ws.Range("AA1").Value = "PinCode [ yes,no ]"
$ws.Range("Z1").Copy($ws.Range("AB1"))
$ws.Range("AB1").Value = "AddressState [ yes,no ]"

# ---------------------------------------------------------------------
# 2. Append "yes" answers for the two new fields on every data row.
# ---------------------------------------------------------------------
$dataRows = @(2, 4, 6, 8, 10, 12, 14)
foreach ($r in $dataRows) {
    $ws.Range("AT$r").Value = "yes"
    $ws.Range("AU$r").Value = "yes"
}

# ---------------------------------------------------------------------
# 3. Row 14 ("Cash"): the collector-name authorisation fields no longer
#    need to be filled in, so flip B/G/H/K from "no" to "yes".
# ---------------------------------------------------------------------
$ws.Range("B14").Value = "yes"
$ws.Range("G14").Value = "yes"
$ws.Range("H14").Value = "yes"
$ws.Range("K14").Value = "yes"

# ---------------------------------------------------------------------
# 4. Column widths for the two freshly inserted header columns.
# ---------------------------------------------------------------------
$ws.Columns("AA").ColumnWidth = 15.8333333333333
$ws.Columns("AB").ColumnWidth = 11.6666666666667

# ---------------------------------------------------------------------
# 5. Restore the workbook's on-screen selection/scroll position.
# ---------------------------------------------------------------------
$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
